# Applies the dual-variable updates for iteration 5 of the branch-and-price
# run with L-shaped subproblems.

$wb = $excel.ActiveWorkbook

# --- Sheet "u_MAB" ---
$wsMAB = $wb.Worksheets.Item("u_MAB")
$wsMAB.Range("B15").Value = 0
$wsMAB.Range("A16").Value = 0
$wsMAB.Range("B40").Value = 0
$wsMAB.Range("B51").Value = 3.088175337449572
$wsMAB.Range("A53").Value = 0.05308741666666674
$wsMAB.Range("A61").Value = 1.419525460645199
$wsMAB.Range("B61").Value = 0.2296470621116087

# --- Sheet "u_EOH" ---
$wsEOH = $wb.Worksheets.Item("u_EOH")
$wsEOH.Range("A3").Value = -3.219476822452852

# --- Sheet "v_l" ---
$wsVL = $wb.Worksheets.Item("v_l")
$wsVL.Range("A2").Value = 5529078.633994577
$wsVL.Range("A3").Value = 5451788.925980421
$wsVL.Range("A4").Value = 3052068.134159257
